$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 2.238867881535716
$ws.Range("D2").Value = 0.03675324539971693

$ws.Range("C4").Value = 0.4485879127129438
$ws.Range("D4").Value = 0.01518174086536896

$ws.Range("C5").Value = 0.2118048282336533
$ws.Range("D5").Value = 0.01902699034649488

$ws.Range("C6").Value = 0.1950786199870751
$ws.Range("D6").Value = 0.01261073097415048

$ws.Range("C8").Value = 0.1409375134690638
$ws.Range("D8").Value = 0.03197758134586337
